$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so that values such as
# "589.64" or "0.0373" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 43 becomes VeChain, Row 44 becomes Maker (rows swapped identity with new values)
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0373"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.928.79"
$ws.Range("E44").Value = "  +0.42%  "


# Remaining D/E numeric and percentage updates
$ws.Range("D2").Value = "63.965.68"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "3.166.96"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "589.64"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "146.59"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.166.15"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +6.92%  "
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "37.18"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "3.668.82"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "63.749.82"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "3.145.45"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "7.13"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "467.29"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "14.34"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "7.50"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "13.08"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("D25").Value = "81.50"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "9.14"
$ws.Range("E28").Value = "  +6.76%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "7.03"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").Value = "27.14"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").Value = "3.39"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").Value = "2.32"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "6.04"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "50.61"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "443.83"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "8.74"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D45").Value = "0.275"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").Value = "36.20"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D48").Value = "125.60"
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "0.111"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "24.51"
$ws.Range("E51").Value = "  -0.61%  "

# Restore the original (default) cell style on column D now that the
# values have been written as text, so no visible style change remains.
$ws.Range("D2:D51").Style = "Normal"
